# Applies the cryptos price/volume refresh described in the commit "Updated cryptos list".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.211.49"
$ws.Range("E2").Value = "  -2.72%  "

$ws.Range("D3").Value = "3.815.64"
$ws.Range("E3").Value = "  +1.92%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'595.92"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -3.83%  "

$ws.Range("D6").Value = "'172.89"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.09%  "

$ws.Range("D7").Value = "3.811.61"
$ws.Range("E7").Value = "  +1.90%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").Value = "'0.535"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.17%  "

$ws.Range("D10").Value = "'0.160"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.75%  "

$ws.Range("D11").Value = "'6.34"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.27%  "

$ws.Range("D12").Value = "'0.471"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.22%  "

$ws.Range("D13").Value = "'38.49"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -4.01%  "

$ws.Range("D14").Value = "'0.0000245"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.72%  "

$ws.Range("D15").Value = "4.437.04"
$ws.Range("E15").Value = "  +1.75%  "

$ws.Range("D16").Value = "3.799.22"
$ws.Range("E16").Value = "  +1.76%  "

$ws.Range("D17").Value = "68.326.95"
$ws.Range("E17").Value = "  -2.54%  "

$ws.Range("E18").Value = "  -4.40%  "

$ws.Range("D19").Value = "'7.26"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -4.14%  "

$ws.Range("D20").Value = "'16.05"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.01%  "

$ws.Range("D21").Value = "'490.26"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.94%  "

$ws.Range("D22").Value = "'9.39"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.82%  "

$ws.Range("D23").Value = "'0.740"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.50%  "

$ws.Range("D24").Value = "'86.07"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.03%  "

$ws.Range("D25").Value = "'2.38"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -7.12%  "

$ws.Range("D26").Value = "'0.0000139"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +5.40%  "

$ws.Range("D27").Value = "'12.28"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -5.47%  "

$ws.Range("D28").Value = "'10.18"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -8.71%  "

$ws.Range("E29").Value = "  +0.15%  "

$ws.Range("D30").Value = "'2.94"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.62%  "

$ws.Range("D31").Value = "'2.44"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.16%  "

$ws.Range("D32").Value = "'33.08"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +8.51%  "

$ws.Range("D33").Value = "'7.65"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.98%  "

$ws.Range("D34").Value = "'0.111"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.25%  "

$ws.Range("E35").Value = "  -0.10%  "

$ws.Range("D36").Value = "'1.00"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -5.03%  "

$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D37").Value = "'5.84"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.03%  "

$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.137"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.48%  "

$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").Value = "'458.47"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +7.20%  "

$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").Value = "'0.326"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -6.58%  "

$ws.Range("D41").Value = "'49.22"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.29%  "

$ws.Range("D42").Value = "'2.01"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.94%  "

$ws.Range("D43").Value = "'2.90"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -8.84%  "

$ws.Range("D44").Value = "'8.35"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.86%  "

$ws.Range("D45").Value = "'41.50"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -6.84%  "

$ws.Range("D46").Value = "2.846.56"
$ws.Range("E46").Value = "  -4.05%  "

$ws.Range("D48").Value = "'0.0353"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.49%  "

$ws.Range("E49").Value = "  +0.62%  "

$ws.Range("D50").Value = "'26.50"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.68%  "

$ws.Range("D51").Value = "'23.53"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +7.94%  "
